$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row content: a label for the new Class Diagram entry
$ws.Range("C7").Value = "Class Diagram for Pedestrian And Car TLS"

# Update the active selection/cell to match the new cursor position
$ws.Range("E9").Select()
